# Generate Report for Handback
# Updates the localization-status workbook to reflect that the two
# outstanding files have been handed back (are now in sync with en-US):
#   - Status column changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" on the Overview sheet and on each
#     language sheet.
#   - Each language sheet gains a "Latest Target File" (E) and
#     "Latest Handback File" (F) hyperlink, mirroring the source file (A)
#     and handoff file (C) respectively.
#   - The "Latest Handback DateTime" (G) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: rows 2 and 3 (columns B and C) share the status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2,2).Value = $newStatus
$wsOverview.Cells.Item(2,3).Value = $newStatus
$wsOverview.Cells.Item(3,2).Value = $newStatus
$wsOverview.Cells.Item(3,3).Value = $newStatus

# --- Language sheets ---
$languages = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-03-09 01:33:03" },
    @{ Name = "de-de"; HandbackTime = "2016-03-09 01:33:37" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Collect existing hyperlink addresses for A2, C2, A3, C3 so the new
    # Target File / Handback File hyperlinks can point at the same targets.
    $urlA2 = ""
    $urlC2 = ""
    $urlA3 = ""
    $urlC3 = ""
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 1) { $urlA2 = $hl.Address }
        if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 3) { $urlC2 = $hl.Address }
        if ($hl.Range.Row -eq 3 -and $hl.Range.Column -eq 1) { $urlA3 = $hl.Address }
        if ($hl.Range.Row -eq 3 -and $hl.Range.Column -eq 3) { $urlC3 = $hl.Address }
    }

    $displayA2 = $ws.Cells.Item(2,1).Text
    $displayC2 = $ws.Cells.Item(2,3).Text
    $displayA3 = $ws.Cells.Item(3,1).Text
    $displayC3 = $ws.Cells.Item(3,3).Text

    # Status text (Ready for handoff -> Handed back: in sync with en-US)
    $ws.Cells.Item(2,2).Value = $newStatus
    $ws.Cells.Item(3,2).Value = $newStatus

    # Latest Target File (E) / Latest Handback File (F) hyperlinks
    $ws.Hyperlinks.Add($ws.Cells.Item(2,5), $urlA2, "", "", $displayA2) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(2,6), $urlC2, "", "", $displayC2) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(3,5), $urlA3, "", "", $displayA3) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item(3,6), $urlC3, "", "", $displayC3) | Out-Null

    # Latest Handback DateTime (G) for the two handed-back rows
    $ws.Cells.Item(2,7).Value = $lang.HandbackTime
    $ws.Cells.Item(3,7).Value = $lang.HandbackTime
}
